# Updates cryptos list values per commit: "Updated cryptos list on Mon Mar  6 03:47:39 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "1.563.34", "1.000") that must
# stay as literal text. Temporarily force Text format on the whole Price column so Excel
# does not reinterpret these strings as numbers, then restore the default style afterward
# so no stray cell styles are introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "22.414.61"
$ws.Range("D3").Value = "1.563.34"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "286.87"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").Value = "0.3652"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").Value = "49.72"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").Value = "1.128"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "0.07405"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "6.870"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "1.562.55"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "89.05"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "0.06730"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D21").Value = "6.294"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").Value = "22.414.36"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "2.550"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "149.36"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "19.65"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").Value = "5.000"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "123.23"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").Value = "1.739.27"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "1.055"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").Value = "6.094"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "9.617"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").Value = "0.08263"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("D37").Value = "0.02391"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("D38").Value = "1.310"
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06377"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2220"
$ws.Range("D41").Value = "5.324"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").Value = "11.20"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "0.6086"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "3.767"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").Value = "0.5742"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").Value = "2.014"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "124.58"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("D50").Value = "1.226"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "0.07246"
$ws.Range("E51").Value = "  -1.13%  "

$priceRange.Style = "Normal"
